$d = $word.ActiveDocument
$t = $d.Tables.Item(1)

function Set-CellText($table, $row, $col, $newText) {
    $cell = $table.Cell($row, $col)
    $cell.Range.Text = $newText
}

# Row 2 (ММП v2 block, СКО=0.0)
Set-CellText $t 2 1 "ММП"
Set-CellText $t 2 12 "0.0038"
Set-CellText $t 2 13 "0.0575"

# Row 3 (Метод N пеленгов, СКО=0.0)
Set-CellText $t 3 13 "0.0022"

# Row 4 (ММП v2 block, СКО=0.1)
Set-CellText $t 4 1 "ММП"
Set-CellText $t 4 12 "0.0028"
Set-CellText $t 4 13 "0.0139"

# Row 5 (Метод N пеленгов, СКО=0.1)
Set-CellText $t 5 13 "0.0013"

# Row 6 (ММП v2 block, СКО=0.2)
Set-CellText $t 6 1 "ММП"
Set-CellText $t 6 12 "0.0031"
Set-CellText $t 6 13 "0.0152"

# Row 7 (Метод N пеленгов, СКО=0.2)
Set-CellText $t 7 12 "0.0003"
Set-CellText $t 7 13 "0.0028"

# Row 8 (ММП v2 block, СКО=0.3)
Set-CellText $t 8 1 "ММП"
Set-CellText $t 8 12 "0.0038"
Set-CellText $t 8 13 "0.0664"

# Row 9 (Метод N пеленгов, СКО=0.3)
Set-CellText $t 9 13 "0.0050"

# Row 10 (ММП v2 block, СКО=0.5)
Set-CellText $t 10 1 "ММП"
Set-CellText $t 10 12 "0.0036"
Set-CellText $t 10 13 "0.0427"

# Row 11 (Метод N пеленгов, СКО=0.5)
Set-CellText $t 11 13 "0.0010"

# Row 12 (ММП v2 block, СКО=1.0)
Set-CellText $t 12 1 "ММП"
Set-CellText $t 12 12 "0.0040"
Set-CellText $t 12 13 "0.0603"

# Row 13 (Метод N пеленгов, СКО=1.0)
Set-CellText $t 13 13 "0.0082"
